$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(100, 101, 102, 103, 104, 105, 106, 107, 114, 115)
$newValues = @{
    100 = 188075.5791
    101 = 341996.99
    102 = 1726834.841
    103 = 1100083.853
    104 = 144673.5224
    105 = 204738.3449
    106 = 189951.8325
    107 = 1011890.156
    114 = 9352.117878999999
    115 = 726146.4398000001
}

foreach ($r in $rows) {
    $ws.Range("J" + $r + ":AS" + $r).Value = $newValues[$r]
}
